# Refresh the config sheet so BANK/FOLDER values point at today's run and
# STATEMENT_DATE uses today's date instead of the stale one.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Drop the leftover explicit "default" formatting on the label/explanation
# columns (and the stray empty B3 cell) so the sheet goes back to plain
# General formatting everywhere except the date cell below.
$clearFormatCells = @("A2","C2","A3","B3","C3","A4","C4","A5","C5","A6","C6","A7","C7")
foreach ($addr in $clearFormatCells) {
    $ws.Range($addr).ClearFormats()
}

# STATEMENT_DATE -> today's date (kept as text, matching the original format).
$ws.Range("B8").Value = "27 OCT 2022"
$ws.Range("B8").NumberFormat = "@"

# Point the bot's source/destination reconciliation report + folder path at
# today's files.
$ws.Range("B4").Value = "NOSTRO RECO 26TH OCT 2022.xlsx"
$ws.Range("B7").Value = "C:\Users\RPA\Desktop\TestFolder\"
$ws.Range("B5").Value = "Nostro Reco Report 27th Oct 2022_.xlsx"

# B3 (CITRIX_FILE_NAME value) was only ever an empty placeholder cell - drop it.
$ws.Range("B3").ClearContents()

$ws.Range("B5").Select() | Out-Null
